# Applies the scheduled-runner market data refresh described in the commit
# message "chore: update Sheets via scheduled runner" to Sheets/Ravana_Profits.xlsx.
# Each per-job worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) stores cached
# market-board pricing/profit figures (columns H:N) per leve row; this script
# rewrites the specific cells whose cached values changed between pulls.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 11997.5
$ws.Range("J32").Value = 12663.333
$ws.Range("L32").Value = 12663.333
$ws.Range("N32").Value = -13315.333
$ws.Range("H43").Value = 1766.6666
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H116").Value = 4044.7144
$ws.Range("I116").Value = 3215
$ws.Range("J116").Value = 4874.4287
$ws.Range("K116").Value = 3215
$ws.Range("L116").Value = 4874.4287
$ws.Range("M116").Value = 227
$ws.Range("N116").Value = -11758.4287
$ws.Range("H132").Value = 1735.4166
$ws.Range("I132").Value = 1583.6
$ws.Range("K132").Value = 4750.799999999999
$ws.Range("M132").Value = -2220.799999999999
$ws.Range("H137").Value = 2134.2856
$ws.Range("I137").Value = 848.3333
$ws.Range("J137").Value = 3098.75
$ws.Range("K137").Value = 2544.9999
$ws.Range("L137").Value = 9296.25
$ws.Range("M137").Value = 5.000100000000202
$ws.Range("N137").Value = -14396.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1065.8667
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H32").Value = 6453.8
$ws.Range("I32").Value = 5893.1665
$ws.Range("K32").Value = 5893.1665
$ws.Range("M32").Value = -5606.1665
$ws.Range("H61").Value = 1657.7
$ws.Range("I61").Value = 716
$ws.Range("K61").Value = 716
$ws.Range("M61").Value = -504
$ws.Range("H74").Value = 13329831
$ws.Range("I74").Value = 18175770
$ws.Range("K74").Value = 18175770
$ws.Range("M74").Value = -18174896
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H77").Value = 13329831
$ws.Range("I77").Value = 18175770
$ws.Range("K77").Value = 90878850
$ws.Range("M77").Value = -90874482
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H116").Value = 1065.8667
$ws.Range("I116").Value = 1000
$ws.Range("K116").Value = 1000
$ws.Range("M116").Value = 1294
$ws.Range("H136").Value = 1657.7
$ws.Range("I136").Value = 716
$ws.Range("K136").Value = 2148
$ws.Range("M136").Value = 402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1065.8667
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -886
$ws.Range("H64").Value = 1333.3334
$ws.Range("J64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("N64").Value = -1950
$ws.Range("H67").Value = 1333.3334
$ws.Range("J67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("N67").Value = -3060
$ws.Range("H75").Value = 10000
$ws.Range("I75").Value = 10000
$ws.Range("K75").Value = 10000
$ws.Range("M75").Value = -9064
$ws.Range("H78").Value = 10000
$ws.Range("I78").Value = 10000
$ws.Range("K78").Value = 30000
$ws.Range("M78").Value = -25320
$ws.Range("H134").Value = 2520
$ws.Range("I134").Value = 2308.5715
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6925.7145
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4390.7145
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2845.6667
$ws.Range("I31").Value = 3012
$ws.Range("J31").Value = 2014
$ws.Range("K31").Value = 3012
$ws.Range("L31").Value = 2014
$ws.Range("M31").Value = -2717
$ws.Range("N31").Value = -2604
$ws.Range("H34").Value = 2845.6667
$ws.Range("I34").Value = 3012
$ws.Range("J34").Value = 2014
$ws.Range("K34").Value = 3012
$ws.Range("L34").Value = 2014
$ws.Range("M34").Value = -2810
$ws.Range("N34").Value = -2418
$ws.Range("H86").Value = 23379.812
$ws.Range("I86").Value = 8774.799999999999
$ws.Range("J86").Value = 47721.5
$ws.Range("K86").Value = 8774.799999999999
$ws.Range("L86").Value = 47721.5
$ws.Range("M86").Value = -7651.799999999999
$ws.Range("N86").Value = -49967.5
$ws.Range("H89").Value = 23379.812
$ws.Range("I89").Value = 8774.799999999999
$ws.Range("J89").Value = 47721.5
$ws.Range("K89").Value = 43874
$ws.Range("L89").Value = 238607.5
$ws.Range("M89").Value = -38258
$ws.Range("N89").Value = -249839.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3841.8572
$ws.Range("I80").Value = 2771.75
$ws.Range("K80").Value = 2771.75
$ws.Range("M80").Value = -1773.75
$ws.Range("H83").Value = 3841.8572
$ws.Range("I83").Value = 2771.75
$ws.Range("K83").Value = 13858.75
$ws.Range("M83").Value = -8866.75
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H102").Value = 2915.625
$ws.Range("I102").Value = 2302.2
$ws.Range("J102").Value = 3938
$ws.Range("K102").Value = 2302.2
$ws.Range("L102").Value = 3938
$ws.Range("M102").Value = -680.1999999999998
$ws.Range("N102").Value = -7182
$ws.Range("H113").Value = 2384.1428
$ws.Range("I113").Value = 2237.8
$ws.Range("K113").Value = 2237.8
$ws.Range("M113").Value = -67.80000000000018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6874.5
$ws.Range("J122").Value = 6833
$ws.Range("L122").Value = 20499
$ws.Range("N122").Value = -25399

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15336
$ws.Range("H122").Value = 2042.3334
$ws.Range("I122").Value = 2047.625
$ws.Range("K122").Value = 6142.875
$ws.Range("M122").Value = -3692.875
$ws.Range("H126").Value = 1275.0667
$ws.Range("I126").Value = 1287.5714
$ws.Range("K126").Value = 3862.7142
$ws.Range("M126").Value = -1392.7142
$ws.Range("H132").Value = 4233.1177
$ws.Range("I132").Value = 4056.8
$ws.Range("J132").Value = 4306.5835
$ws.Range("K132").Value = 12170.4
$ws.Range("L132").Value = 12919.7505
$ws.Range("M132").Value = -9640.400000000001
$ws.Range("N132").Value = -17979.7505
$ws.Range("H136").Value = 2969
$ws.Range("I136").Value = 2750
$ws.Range("K136").Value = 8250
$ws.Range("M136").Value = -5700
